$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(4, 7, "i3-1313", 3.4, 150),
    @(5, 7, "i3-1413", 3.3, 180),
    @(6, 7, "i5-5010", 3.1, 250),
    @(7, 7, "i5-6060", 3.2, 260),
    @(8, 7, "i7-7056", 3.5, 300),
    @(9, 7, "i7-7100", 3.6, 320),
    @(10, 6, "fx-4100", 2.9, 150),
    @(11, 6, "fx-4200", 3.1, 165),
    @(12, 6, "fx-6100", 3.3, 180),
    @(13, 6, "fx-6300", 3.5, 200),
    @(14, 6, "fx-8300 ", 3.8, 250),
    @(15, 6, "fx-9370", 4.2, 300)
)

$r = 5
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r = $r + 1
}
